$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- I7: hours worked on 2/7 entry updated from 0.6 to 0.75 ---
$ws.Range("I7").Value = 0.75

# --- Row 21: Feb 7 team meeting entry ---
$ws.Range("A21").Value = 41677
$ws.Range("A21").NumberFormat = "m/d/yy"
$ws.Range("B21").Value = "Team meeting. Recorded meeting minutes. Continued drafting software requirements specification document. Discussed high-level design problems and solutions. Developed additional client questions."
$ws.Range("I21").Value = 3

# --- Row 22: Feb 7 SRS update entry ---
$ws.Range("A22").Value = 41677
$ws.Range("A22").NumberFormat = "m/d/yy"
$ws.Range("B22").Value = "Heavily updated software requirements specification document by adding sections and updating information."
$ws.Range("I22").Value = 2.25

# --- Row 23: Feb 8 team meeting entry ---
$ws.Range("A23").Value = 41678
$ws.Range("A23").NumberFormat = "m/d/yy"
$ws.Range("B23").Value = "Team meeting. Recorded meeting minutes. Continued drafting software requirements specification document. Contributed to prototype user interface color scheme discussions."
$ws.Range("I23").Value = 4

# --- Row 24: Feb 8 SRS revision completed entry (bottom border of the block) ---
$ws.Range("A24").Value = 41678
$ws.Range("A24").Borders.Item(9).LineStyle = 1
$ws.Range("A24").NumberFormat = "mm-dd-yy"
$ws.Range("B24").Value = "Completed the first major draft revision of the software requirements specification document."
$ws.Range("I24").Value = 1.25

# --- Update the window scroll position / selection to reflect the newly-filled rows ---
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B24:H24").Select()
